# readout.pptx maintenance edit:
#   1. Bump the cached "fixed" footer date from 23.08.2012 to 25.08.2012
#      everywhere it is cached (slide master, all 11 slide layouts, notes
#      master).
#   2. Nudge two textbox labels ("measurement" / "latching") on slide 2
#      sideways (x offset only; y/width/height untouched).

$p = $ppt.ActivePresentation

$oldDate = "23.08.2012"
$newDate = "25.08.2012"

# --- 1a. Slide master + every slide layout: the date lives in a plain
#     shape (placeholder) whose TextFrame we can set directly. ---
function Update-DateInShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateInShapes $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateInShapes $layout.Shapes
}

# --- 1b. Notes master: direct shape text assignment is a no-op for this
#     container in this host, but the HeadersFooters.DateAndTime façade
#     does persist the change, so use that instead. (Note: reading
#     .Text back through this façade yields an empty string even though
#     the write lands correctly, so we set unconditionally rather than
#     gating on the current value.) ---
$notesMaster = $p.NotesMaster
$notesDateAndTime = $notesMaster.HeadersFooters.DateAndTime
$notesDateAndTime.Text = $newDate

# --- 2. Slide 2: shift "Textfeld 57" (measurement) and "Textfeld 58"
#     (latching) horizontally. Shape.Left/.Top are single-precision
#     points under the hood, so converting EMU -> points and back can
#     truncate one EMU low; a hair of epsilon keeps the round-trip exact.
$slide2 = $p.Slides.Item(2)
$emuPerPoint = 12700
$epsilon = 0.000001

$measurementShape = $slide2.Shapes.Item("Textfeld 57")
$measurementShape.Left = (251520 / $emuPerPoint) + $epsilon

$latchingShape = $slide2.Shapes.Item("Textfeld 58")
$latchingShape.Left = (1297732 / $emuPerPoint) + $epsilon
